$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")
$c = $ws.Range("G2")
$c.Value = 1
$c.Borders.LineStyle = 1
